$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F5 to new placeholder text "<<jumlah penghuni>>"
# (leading apostrophe forces text/quote-prefix entry, preserving the
# cell's existing quotePrefix style instead of Excel re-normalizing it)
$ws.Range("F5").Value = "'<<jumlah penghuni>>"

# Update the active selection to G7
$ws.Range("G7").Select()
